$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $s = $ws.Range($cellRef).Style
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = $s
}

Set-TextValue 'D2' '55.331.94'
Set-TextValue 'E2' '  +1.59%  '
Set-TextValue 'D3' '2.300.06'
Set-TextValue 'E3' '  +0.56%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.17%  '
Set-TextValue 'D5' '508.08'
Set-TextValue 'E5' '  +1.08%  '
Set-TextValue 'D6' '130.06'
Set-TextValue 'E6' '  -0.16%  '
Set-TextValue 'D7' '0.995'
Set-TextValue 'E7' '  -0.28%  '
Set-TextValue 'D8' '0.531'
Set-TextValue 'E8' '  +0.30%  '
Set-TextValue 'D9' '2.326.95'
Set-TextValue 'E9' '  +1.18%  '
Set-TextValue 'D10' '0.0983'
Set-TextValue 'E10' '  +2.47%  '
Set-TextValue 'E11' '  +1.78%  '
Set-TextValue 'D12' '5.10'
Set-TextValue 'E12' '  +7.94%  '
Set-TextValue 'D13' '0.342'
Set-TextValue 'E13' '  +1.53%  '
Set-TextValue 'D14' '24.01'
Set-TextValue 'E14' '  +4.36%  '
Set-TextValue 'D15' '2.711.48'
Set-TextValue 'E15' '  +0.56%  '
Set-TextValue 'D16' '55.113.09'
Set-TextValue 'E16' '  +1.25%  '
Set-TextValue 'E17' '  +1.48%  '
Set-TextValue 'D18' '2.315.62'
Set-TextValue 'E18' '  +0.03%  '
Set-TextValue 'D19' '10.76'
Set-TextValue 'E19' '  +4.41%  '
Set-TextValue 'D20' '4.20'
Set-TextValue 'E20' '  +0.82%  '
Set-TextValue 'E21' '  +4.32%  '
Set-TextValue 'D22' '311.69'
Set-TextValue 'E22' '  +2.27%  '
Set-TextValue 'D23' '0.998'
Set-TextValue 'E23' '  -0.16%  '
Set-TextValue 'D24' '60.40'
Set-TextValue 'E24' '  -2.54%  '
Set-TextValue 'D25' '0.993'
Set-TextValue 'E25' '  -0.62%  '
Set-TextValue 'E26' '  -0.02%  '
Set-TextValue 'D27' '7.54'
Set-TextValue 'E27' '  +2.54%  '
Set-TextValue 'D28' '172.77'
Set-TextValue 'E28' '  -0.21%  '
Set-TextValue 'E29' '  +2.79%  '
Set-TextValue 'E30' '  +2.59%  '
Set-TextValue 'B31' 'PancakeSwap'
Set-TextValue 'C31' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D31' '1.64'
Set-TextValue 'E31' '  +0.48%  '
Set-TextValue 'B32' 'Fetch.AI'
Set-TextValue 'C32' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D32' '1.15'
Set-TextValue 'E32' '  +4.41%  '
Set-TextValue 'D33' '18.13'
Set-TextValue 'E33' '  +1.33%  '
Set-TextValue 'D34' '0.999'
Set-TextValue 'E34' '  +0.00%  '
Set-TextValue 'D35' '0.997'
Set-TextValue 'E35' '  +0.04%  '
Set-TextValue 'D36' '1.24'
Set-TextValue 'E36' '  +2.59%  '
Set-TextValue 'D37' '0.919'
Set-TextValue 'E37' '  -5.36%  '
Set-TextValue 'D38' '3.92'
Set-TextValue 'E38' '  +3.50%  '
Set-TextValue 'D39' '36.83'
Set-TextValue 'E39' '  +2.20%  '
Set-TextValue 'E40' '  +1.90%  '
Set-TextValue 'E41' '  +1.03%  '
Set-TextValue 'D42' '136.09'
Set-TextValue 'E42' '  +8.33%  '
Set-TextValue 'E43' '  +1.32%  '
Set-TextValue 'D44' '4.95'
Set-TextValue 'E44' '  +1.02%  '
Set-TextValue 'D45' '260.91'
Set-TextValue 'E45' '  +7.01%  '
Set-TextValue 'D46' '0.0506'
Set-TextValue 'E46' '  +1.60%  '
Set-TextValue 'E47' '  +1.95%  '
Set-TextValue 'D48' '0.555'
Set-TextValue 'E48' '  +0.80%  '
Set-TextValue 'D49' '0.379'
Set-TextValue 'E49' '  +1.33%  '
Set-TextValue 'D50' '0.0211'
Set-TextValue 'E50' '  +1.93%  '
Set-TextValue 'D51' '10.82'
Set-TextValue 'E51' '  +0.36%  '
